{"js": "// Update each two-digit-number \u00f7 one-digit-number expression in the\n// practice-sheet table to the newly generated problem/answer prompt.\n// Every \"old=>new\" pair below is unique in the document, so a scoped\n// search-and-replace per pair is unambiguous and order-independent.\nconst replacements = [\n  [\"40\u00f73=\", \"32\u00f76=\"],\n  [\"92\u00f73=\", \"71\u00f77=\"],\n  [\"15\u00f77=\", \"61\u00f78=\"],\n  [\"20\u00f75=\", \"25\u00f75=\"],\n  [\"68\u00f78=\", \"70\u00f72=\"],\n  [\"46\u00f77=\", \"58\u00f74=\"],\n  [\"57\u00f74=\", \"49\u00f76=\"],\n  [\"61\u00f79=\", \"24\u00f77=\"],\n  [\"18\u00f79=\", \"92\u00f78=\"],\n  [\"56\u00f79=\", \"60\u00f76=\"],\n  [\"13\u00f72=\", \"94\u00f75=\"],\n  [\"91\u00f78=\", \"57\u00f73=\"],\n  [\"62\u00f77=\", \"80\u00f72=\"],\n  [\"95\u00f75=\", \"60\u00f76=\"],\n  [\"34\u00f78=\", \"37\u00f75=\"],\n  [\"28\u00f74=\", \"58\u00f74=\"],\n  [\"86\u00f72=\", \"20\u00f77=\"],\n  [\"40\u00f72=\", \"77\u00f73=\"],\n  [\"98\u00f72=\", \"30\u00f79=\"],\n  [\"84\u00f75=\", \"57\u00f78=\"],\n  [\"36\u00f76=\", \"74\u00f73=\"],\n  [\"80\u00f76=\", \"31\u00f75=\"],\n  [\"45\u00f79=\", \"55\u00f74=\"],\n  [\"68\u00f74=\", \"62\u00f76=\"],\n  [\"42\u00f73=\", \"31\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update each two-digit-number / one-digit-number division expression\n# in the practice-sheet table to the newly generated problem prompt.\n# Every \"old => new\" pair is unique in the document, so a single\n# Find/Replace (wdReplaceAll = 2) per pair is unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"40\u00f73=\", \"32\u00f76=\"),\n    @(\"92\u00f73=\", \"71\u00f77=\"),\n    @(\"15\u00f77=\", \"61\u00f78=\"),\n    @(\"20\u00f75=\", \"25\u00f75=\"),\n    @(\"68\u00f78=\", \"70\u00f72=\"),\n    @(\"46\u00f77=\", \"58\u00f74=\"),\n    @(\"57\u00f74=\", \"49\u00f76=\"),\n    @(\"61\u00f79=\", \"24\u00f77=\"),\n    @(\"18\u00f79=\", \"92\u00f78=\"),\n    @(\"56\u00f79=\", \"60\u00f76=\"),\n    @(\"13\u00f72=\", \"94\u00f75=\"),\n    @(\"91\u00f78=\", \"57\u00f73=\"),\n    @(\"62\u00f77=\", \"80\u00f72=\"),\n    @(\"95\u00f75=\", \"60\u00f76=\"),\n    @(\"34\u00f78=\", \"37\u00f75=\"),\n    @(\"28\u00f74=\", \"58\u00f74=\"),\n    @(\"86\u00f72=\", \"20\u00f77=\"),\n    @(\"40\u00f72=\", \"77\u00f73=\"),\n    @(\"98\u00f72=\", \"30\u00f79=\"),\n    @(\"84\u00f75=\", \"57\u00f78=\"),\n    @(\"36\u00f76=\", \"74\u00f73=\"),\n    @(\"80\u00f76=\", \"31\u00f75=\"),\n    @(\"45\u00f79=\", \"55\u00f74=\"),\n    @(\"68\u00f74=\", \"62\u00f76=\"),\n    @(\"42\u00f73=\", \"31\u00f73=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # wdFindContinue=1, wdReplaceAll=2\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
